$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.938.68'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '2.929.85'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '356.01'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '110.98'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('E7').Value = '  +1.86%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.627'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.47'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0880'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.136'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.90'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').Value = '3.392.05'
$ws.Range('E15').Value = '  +1.29%  '
$ws.Range('D16').Value = '2.918.66'
$ws.Range('E16').Value = '  +1.36%  '
$ws.Range('E17').Value = '  -1.70%  '
$ws.Range('D18').Value = '51.933.80'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.30'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.57'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.02'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.38%  '
$ws.Range('D22').Value = '0.0₃0982'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.98'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '270.95'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.83'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.185'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +12.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.20'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.40'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +14.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.106'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +12.27%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '10.61'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.76%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '38.97'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.06'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '52.08'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0445'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.97%  '
$ws.Range('B36').Value = 'Toncoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.94'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -13.34%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  -1.84%  '
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.75'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.31%  '
$ws.Range('E42').Value = '  +2.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.19'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '119.42'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.17'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.74%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.47'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.45%  '
$ws.Range('D48').Value = '2.140.29'
$ws.Range('E48').Value = '  -2.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.249'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -8.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0334'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.84%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '62.18'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.51%  '
